# "Generate Report for Archive" - refresh the localization-status report:
#   1. The handoff status text moves from "Ready for handoff" to
#      "In Translation" everywhere it's shown (Overview summary columns
#      for zh-cn/de-de, plus the per-language Status column).
#   2. Those status columns get narrower now that the new status text is
#      shorter than the old one, so re-autofit them to the new content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) show the status ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value2 = $newStatus
    }
}
# Narrow columns E and F to match the shorter status text.
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

# --- Per-language sheets: column C is the Status column ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 13.4101845877511
}
